# Auto update Excel log
# Append new PRESENCE_DETECTED rows to the "mmWave" sheet (rows 24-29),
# mirroring the existing log rows (Date / Timestamp / Hour / Location / Value / Status).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$newRows = @(
    @("2026-02-01", "15:55:25", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "15:55:35", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "15:55:46", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "15:55:56", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "15:56:07", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "15:56:17", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 24
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Columns A and C hold date/time-looking text ("2026-02-01", "15:00") that must
    # stay as literal text (matching the rest of the log) instead of being
    # auto-converted into date/time serial numbers.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]

    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[2]

    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
